$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.358.70"
$ws.Range("E2").Value = "  +0.37%  "

$ws.Range("D3").Value = "3.593.55"
$ws.Range("E3").Value = "  +0.23%  "

$ws.Range("E4").Value = "  +0.27%  "

$ws.Range("D5").Value = "'200.13"
$ws.Range("E5").Value = "  +1.99%  "

$ws.Range("D6").Value = "'593.59"
$ws.Range("E6").Value = "  -1.98%  "

$ws.Range("D7").Value = "'0.628"
$ws.Range("E7").Value = "  +0.56%  "

$ws.Range("E9").Value = "  +5.80%  "

$ws.Range("D10").Value = "'0.644"
$ws.Range("E10").Value = "  -1.09%  "

$ws.Range("D11").Value = "'53.34"
$ws.Range("E11").Value = "  -1.00%  "

$ws.Range("E12").Value = "  -1.16%  "

$ws.Range("D13").Value = "'9.63"
$ws.Range("E13").Value = "  +0.69%  "

$ws.Range("D14").Value = "'699.34"
$ws.Range("E14").Value = "  +16.95%  "

$ws.Range("D15").Value = "4.167.97"
$ws.Range("E15").Value = "  +1.13%  "

$ws.Range("D16").Value = "70.460.96"
$ws.Range("E16").Value = "  +0.41%  "

$ws.Range("D17").Value = "'12.73"
$ws.Range("E17").Value = "  -1.81%  "

$ws.Range("D18").Value = "'19.04"
$ws.Range("E18").Value = "  -1.29%  "

$ws.Range("D19").Value = "3.594.54"
$ws.Range("E19").Value = "  +1.57%  "

$ws.Range("E20").Value = "  +0.39%  "

$ws.Range("D21").Value = "'0.993"
$ws.Range("E21").Value = "  +0.05%  "

$ws.Range("D22").Value = "'18.20"
$ws.Range("E22").Value = "  +2.03%  "

$ws.Range("D23").Value = "'111.12"
$ws.Range("E23").Value = "  +8.05%  "

$ws.Range("D24").Value = "'5.32"
$ws.Range("E24").Value = "  +3.24%  "

$ws.Range("D25").Value = "'4.53"
$ws.Range("E25").Value = "  -2.34%  "

$ws.Range("D26").Value = "'3.00"
$ws.Range("E26").Value = "  -2.41%  "

$ws.Range("D27").Value = "'10.50"
$ws.Range("E27").Value = "  -3.18%  "

$ws.Range("E28").Value = "  -0.69%  "

$ws.Range("D29").Value = "'9.98"
$ws.Range("E29").Value = "  +4.06%  "

$ws.Range("D30").Value = "'34.69"
$ws.Range("E30").Value = "  +2.96%  "

$ws.Range("D31").Value = "'4.42"
$ws.Range("E31").Value = "  +2.54%  "

$ws.Range("D32").Value = "'7.06"
$ws.Range("E32").Value = "  -0.66%  "

$ws.Range("D33").Value = "'12.24"
$ws.Range("E33").Value = "  -0.92%  "

$ws.Range("E34").Value = "  -1.85%  "

$ws.Range("D35").Value = "'63.56"
$ws.Range("E35").Value = "  +0.22%  "

$ws.Range("D36").Value = "0.0₃0846"
$ws.Range("E36").Value = "  +2.76%  "

$ws.Range("D37").Value = "3.801.66"
$ws.Range("E37").Value = "  -0.86%  "

$ws.Range("D38").Value = "'0.999"
$ws.Range("E38").Value = "  -0.18%  "

$ws.Range("D39").Value = "'3.64"
$ws.Range("E39").Value = "  +0.76%  "

$ws.Range("B40").Value = "Bittensor"
$ws.Range("C40").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D40").Value = "'509.54"
$ws.Range("E40").Value = "  -2.23%  "

$ws.Range("B41").Value = "Fetch.AI"
$ws.Range("C41").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D41").Value = "'3.00"
$ws.Range("E41").Value = "  -7.58%  "

$ws.Range("D42").Value = "'36.42"
$ws.Range("E42").Value = "  -1.04%  "

$ws.Range("D43").Value = "'0.381"
$ws.Range("E43").Value = "  -3.24%  "

$ws.Range("D44").Value = "'0.137"
$ws.Range("E44").Value = "  +2.48%  "

$ws.Range("D45").Value = "'0.0470"
$ws.Range("E45").Value = "  +3.51%  "

$ws.Range("D46").Value = "'3.04"
$ws.Range("E46").Value = "  +7.41%  "

$ws.Range("D47").Value = "'3.42"
$ws.Range("E47").Value = "  +3.12%  "

$ws.Range("E48").Value = "  +0.38%  "

$ws.Range("E49").Value = "  +1.07%  "

$ws.Range("E50").Value = "  -0.15%  "

$ws.Range("D51").Value = "'1.82"
$ws.Range("E51").Value = "  +21.43%  "
